$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-18 (header stays in row 1)
$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $rowIndex++
}

# Remove the old row 19 entirely, since the table now only has 17 data rows (2-18)
$ws.Rows.Item(19).Delete()
